$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "богдана"
$ws.Range("B2").Value = "Відсутній"
